# Applies the "main branch" marker edit plus the trailing shaded paragraph,
# matching the target unified diff.
$d = $word.ActiveDocument

# --- 1. First paragraph: append two trailing spaces, then a red
#        "(This is a change - Version for main branch)" marker, built out
#        of three separate runs (mirroring how Word split them on save). ---
$p1 = $d.Paragraphs(1)
$r1 = $p1.Range
$r1.End = $r1.End - 1
$r1.Text = "This is a Microsoft word document.  "

$enDash = [char]0x2013

$start1 = $r1.End
$r1.Collapse(0)
$r1.InsertAfter("(This is a change " + $enDash + " Ve")
$run1 = $d.Range($start1, $r1.End)
$run1.Font.Color = 255

$start2 = $r1.End
$r1.Collapse(0)
$r1.InsertAfter("rsion for main branch")
$run2 = $d.Range($start2, $r1.End)
$run2.Font.Color = 255

$start3 = $r1.End
$r1.Collapse(0)
$r1.InsertAfter(")")
$run3 = $d.Range($start3, $r1.End)
$run3.Font.Color = 255

# --- 2. Append a brand-new, empty, shaded paragraph at the very end of the
#        document body (right before sectPr). Building it from literal OOXML
#        (rather than InsertParagraphAfter) avoids inheriting the previous
#        paragraph's run/paragraph-mark formatting. ---
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="F9F9F9"/></w:pPr></w:p>')
